# edit.ps1 -- applies the "black holes" -> "Astronomy" rewrite described by
# the target diff, via the Word COM object model.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: insert `$text` right after `$rng` (which is collapsed/duplicated
# internally) using formatting cloned from the author-supplied parameters.
# Word/this engine coalesces an InsertAfter into the previous run when the
# resulting run properties are byte-identical, so we toggle Bold on/off to
# force a genuine run boundary (mirrors the many tiny runs Word itself
# leaves behind from autocorrect/spell-check) and then stamp the desired
# rPr explicitly so the new run's formatting is correct regardless of what
# it inherited.
# ---------------------------------------------------------------------------
function Insert-RunAfter($rng, $text, $fontName, $size, $setSize) {
    $new = $rng.Duplicate
    $new.Collapse(0)
    $new.InsertAfter($text)
    $new.Font.Bold = 1
    $new.Font.Bold = 0
    $new.Font.Name = $fontName
    $new.Font.Color = 0
    if ($setSize) {
        $new.Font.Size = $size
    }
    return $new
}

# ---------------------------------------------------------------------------
# 1. Global font rename: TimesNewToman -> Times New Roman on every paragraph.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $p.Range.Font.Name = "Times New Roman"
}

# ---------------------------------------------------------------------------
# 2. Title
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Unveiling the Enigma of Black Holes", $false, $false, $false, $false, $false, $true, 1, $false, "The Enigmatic Symphony of Astronomy: Unveiling the Cosmos", 2)

# ---------------------------------------------------------------------------
# 3. Author line: "Joanna Carter" -> "Oliver W" + "." + " Florence" (3 runs)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Joanna Carter")
$rng.Text = "Oliver W"
$rng.Font.Name = "Times New Roman"
$rng.Font.Color = 0
$rng.Font.Size = 18
$r = Insert-RunAfter $rng "." "Times New Roman" 18 $true
$r = Insert-RunAfter $r " Florence" "Times New Roman" 18 $true

# ---------------------------------------------------------------------------
# 4. Email line:
#    "jcarter@astroscience" + "." + "edu"
#    -> "oliver" + "." + "florence@eduworld" + "." + "com"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("jcarter@astroscience", $false, $false, $false, $false, $false, $true, 1, $false, "oliver", 2)

$rng = $d.Content
$rng.Find.Execute("edu")
$rng.Text = "florence@eduworld"
$rng.Font.Name = "Times New Roman"
$rng.Font.Color = 0
$rng.Font.Size = 16
$r = Insert-RunAfter $rng "." "Times New Roman" 16 $true
$r = Insert-RunAfter $r "com" "Times New Roman" 16 $true

# ---------------------------------------------------------------------------
# 5. Main body paragraph (font size 24 half-points = 12pt)
# ---------------------------------------------------------------------------
function Replace-One($old, $new) {
    $r = $d.Content
    $r.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-One "Across the vast tapestry of the universe, there exist celestial enigmas that captivate and confound our understanding" "In the boundless expanse of the universe, where stars twinkle like celestial diamonds and galaxies spiral in a cosmic dance, Astronomy embarks on a quest to unravel the enigmas that shroud our cosmos"

Replace-One "Among these cosmic wonders, black holes stand as ultimate expressions of gravity's dominance, defying our current knowledge and challenging our very perception of time and space" "From the fiery birth of stars to the gravitational grip of black holes, this scientific discipline delves into the mysteries of the heavens, beckoning us to explore the extraordinary symphony of the universe"

Replace-One "In this exploration, we delve into the enigmatic realm of black holes, unraveling their properties, exploring their formation, pondering their implications for cosmology, and discussing the ongoing quests to unravel their mysteries through observations and theoretical frameworks" "Astronomy unveils the tapestry of celestial bodies, unraveling the interconnectedness of the cosmos, and igniting imaginations with its profound implications for our place in the universe"

Replace-One "The allure of black holes lies in their extreme nature" "Astronomy's grand narrative begins with the restless curiosity of ancient stargazers, who observed the rhythmic movements of celestial bodies and sought patterns amidst the twinkling expanse"

Replace-One "These cosmic entities arise from the gravitational collapse of massive stars, resulting in singularities where matter is crushed to infinite density and spacetime is warped beyond recognition" "Throughout history, civilizations haveYang Wang the heavens, mapping constellations and charting the courses of celestial bodies"

# This run is immediately followed by a new pair of runs (" Armed with..." + ".")
$rng = $d.Content
$rng.Find.Execute("Their gravitational pull is so intense that nothing, not even light, can escape their clutches, hence their name: 'black holes'")
$rng.Text = " The masterful achievements of astronomers like Galileo, Copernicus, and Kepler laid the cornerstone of our understanding of the solar system, overturning age-old beliefs and ushering in a new era of scientific inquiry"
$rng.Font.Name = "Times New Roman"
$rng.Font.Color = 0
$rng.Font.Size = 12
# advance past this run and the pre-existing "." run that already follows it
$rng.Collapse(0)
$rng.MoveEndUntil(".", 1) | Out-Null
$rng.MoveEnd(1, 1) | Out-Null
$r = Insert-RunAfter $rng " Armed with telescopes that pierce the veil of darkness, astronomers have embarked on an epic voyage of discovery, pushing the boundaries of human knowledge and expanding our perception of reality" "Times New Roman" 12 $true
$r = Insert-RunAfter $r "." "Times New Roman" 12 $true

Replace-One "Despite their elusive nature, astronomers have gathered compelling evidence for the existence of black holes" "The cosmos is a symphony of cosmic phenomena, where stars, planets, galaxies, and nebulae play their part in an eternal cosmic dance"

Replace-One "Through observations of binary star systems, astronomers have detected the presence of invisible objects exerting gravitational forces far beyond what normal stars could account for" "Gravity, the invisible conductor, orchestrates the movements of celestial bodies, guiding them in their celestial ballet"

# This run is immediately followed by a new pair of runs ("." + " From the spectacular...")
$rng = $d.Content
$rng.Find.Execute("These observations, combined with the theoretical predictions of general relativity, strongly suggest the existence of these cosmic behemoths")
$rng.Text = " The birth, life, and death of stars are a testament to the dynamic nature of the universe, a continuous cycle of creation and destruction"
$rng.Font.Name = "Times New Roman"
$rng.Font.Color = 0
$rng.Font.Size = 12
$rng.Collapse(0)
$rng.MoveEndUntil(".", 1) | Out-Null
$rng.MoveEnd(1, 1) | Out-Null
$r = Insert-RunAfter $rng "." "Times New Roman" 12 $true
$r = Insert-RunAfter $r " From the spectacular explosion of a supernova to the relentless pull of a black hole, the cosmos reveals a hidden realm of beauty and chaos, of unimaginable power and awe-inspiring grandeur" "Times New Roman" 12 $true

# ---------------------------------------------------------------------------
# 6. Summary paragraph (no explicit w:sz -> leave font size untouched)
# ---------------------------------------------------------------------------
Replace-One "The existence of black holes, with their immense gravitational pull and enigmatic properties, has profoundly impacted our understanding of the universe" "Astronomy, an enchanting fusion of science and wonder, unveils the mysteries of the universe, captivating the human imagination with its profound revelations"

Replace-One "From their formation through gravitational collapse to their implications for cosmology, black holes continue to captivate scientists and inspire awe in all who contemplate their existence" "From the birth of stars to the enigma of black holes, from the rhythmic dance of planets to the grandeur of distant galaxies, Astronomy provides a lens through which we glimpse the extraordinary symphony of the cosmos"

# The 3rd sentence collapses down to a single space, followed by new runs:
#   <w:lastRenderedPageBreak/> + "Through observation..." + "." + " As we continue..."
$rng = $d.Content
$rng.Find.Execute("As our knowledge expands and technology advances, the quest to unravel the mysteries of black holes remains an ongoing endeavor, promising to reveal even more profound insights into the nature of our universe")
$rng.Text = " "
$rng.Font.Name = "Times New Roman"
$rng.Font.Color = 0
$rng.Collapse(0)

$r = $rng.Duplicate
$r.InsertAfter("Through observation, experimentation, and relentless exploration, astronomers have illuminated the tapestry of the heavens, broadening our understanding of the universe and our place within it")
$r.Font.Bold = 1
$r.Font.Bold = 0
$r.Font.Name = "Times New Roman"
$r.Font.Color = 0
$r.InsertBefore("")
$pageBreakRange = $r.Duplicate
$pageBreakRange.Collapse(1)
$pageBreakRange.InsertBefore("`v")

$r = Insert-RunAfter $r "." "Times New Roman" 0 $false
$r = Insert-RunAfter $r " As we continue to unravel the secrets of the cosmos, we are reminded of the boundless nature of knowledge and the profound interconnectedness of all things" "Times New Roman" 0 $false

# ---------------------------------------------------------------------------
# 7. Trailing empty paragraph at the very end of the document body.
# ---------------------------------------------------------------------------
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

Write-Host "edit complete"
